$wb = $excel.ActiveWorkbook

# Data for the new row 53 on each worksheet: only column A (timestamp) changes
# relative to existing row 52; columns B-I are identical copies of row 52.
$rowData = @(
    @{ Sheet = "MID_LFT_#1"; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x68"; E = "0x07"; F = 400;  G = [double]"5.68631262647113E+23"; H = 360; I = 7  }
    @{ Sheet = "MID_LFT_#2"; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x60"; E = "0x19"; F = 380;  G = [double]"5.68432987514711E+23"; H = 352; I = 25 }
    @{ Sheet = "MID_PLT_#1"; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x69"; E = "0x15"; F = 110;  G = [double]"5.68631262647113E+23"; H = 105; I = 15 }
    @{ Sheet = "MID_PLT_#2"; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x7E"; E = "0x9";  F = 130;  G = [double]"5.68631262647113E+23"; H = 126; I = 9  }
)

$newDate = 45839.46072916667

foreach ($item in $rowData) {
    $ws = $wb.Worksheets.Item($item.Sheet)

    $ws.Cells.Item(53, 1).Value = $newDate
    $ws.Cells.Item(53, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item(53, 2).Value = $item.B
    $ws.Cells.Item(53, 3).Value = $item.C
    $ws.Cells.Item(53, 4).Value = $item.D
    $ws.Cells.Item(53, 5).Value = $item.E
    $ws.Cells.Item(53, 6).Value = $item.F
    $ws.Cells.Item(53, 7).Value = $item.G
    $ws.Cells.Item(53, 8).Value = $item.H
    $ws.Cells.Item(53, 9).Value = $item.I
}
